# Add a new note (with a hyperlink to an uploaded Excel file) to the
# document. The note replaces the first of the two blank paragraphs that
# follow the "Conclusion" paragraph; the final blank paragraph before the
# section properties is left untouched.

$d = $word.ActiveDocument

# Locate the "Conclusion" body paragraph so we anchor off content rather
# than a hard-coded index, then grab the blank paragraph right after it.
$n = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $n; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Future analyses could expand on this work*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not find the anchor paragraph."
}

$target = $d.Paragraphs.Item($anchorIndex + 1)
$targetRange = $target.Range

# Add the hyperlink run first, while the paragraph is still empty - the
# hyperlink is always inserted at the start of the paragraph it is given.
$hyperlink = $d.Hyperlinks.Add(
    $targetRange,
    "https://1drv.ms/x/c/986a3108c35e0a35/EbFJ3example_Crypto_Data_API?e=shared",
    $null,
    $null,
    "Crypto_Data_API.xlsx")

# Re-fetch the paragraph/range and remember the (now fixed) insertion
# point just before the hyperlink run.
$target = $d.Paragraphs.Item($anchorIndex + 1)
$insertPos = $target.Range.Start

$lineBreak = [char]11

# Insert the remaining runs in reverse order, each one immediately before
# the hyperlink, so the final reading order is:
#   run1: plain sentence
#   run2: line break + sentence
#   run3: line break
#   hyperlink: Crypto_Data_API.xlsx

$r3 = $d.Range($insertPos, $insertPos)
$r3.InsertBefore($lineBreak)

$r2 = $d.Range($insertPos, $insertPos)
$r2.InsertBefore("$($lineBreak)I am attaching the link of the file for the same below : ")

$r1 = $d.Range($insertPos, $insertPos)
$r1.InsertBefore("I have also done this task using the excel itself by fetching the data from web using the API URL.")

Write-Output "Inserted note paragraph with hyperlink."
